$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 98.703705
$ws.Range("I2").Value = 98.69231000000001
$ws.Range("J2").Value = 99
$ws.Range("K2").Value = 98.69231000000001
$ws.Range("L2").Value = 99
$ws.Range("M2").Value = 14.30768999999999
$ws.Range("N2").Value = -325
$ws.Range("H9").Value = 89
$ws.Range("I9").Value = 89
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 89
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 80
$ws.Range("N9").ClearContents()
$ws.Range("H40").Value = 1690.5264
$ws.Range("I40").Value = 1280
$ws.Range("J40").Value = 1800
$ws.Range("K40").Value = 1280
$ws.Range("L40").Value = 1800
$ws.Range("M40").Value = -1105
$ws.Range("N40").Value = -2150
$ws.Range("H58").Value = 7981.1055
$ws.Range("I58").Value = 220.08333
$ws.Range("J58").Value = 21285.715
$ws.Range("K58").Value = 660.24999
$ws.Range("L58").Value = 63857.145
$ws.Range("M58").Value = -510.24999
$ws.Range("N58").Value = -64157.145
$ws.Range("H141").Value = 2692.4707
$ws.Range("I141").Value = 1844.5454
$ws.Range("J141").Value = 3098
$ws.Range("K141").Value = 5533.6362
$ws.Range("L141").Value = 9294
$ws.Range("M141").Value = -353.6361999999999
$ws.Range("N141").Value = -19654

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1413.8
$ws.Range("I2").Value = 948.5714
$ws.Range("J2").Value = 2499.3333
$ws.Range("K2").Value = 948.5714
$ws.Range("L2").Value = 2499.3333
$ws.Range("M2").Value = -835.5714
$ws.Range("N2").Value = -2725.3333
$ws.Range("H116").Value = 1413.8
$ws.Range("I116").Value = 948.5714
$ws.Range("J116").Value = 2499.3333
$ws.Range("K116").Value = 948.5714
$ws.Range("L116").Value = 2499.3333
$ws.Range("M116").Value = 1345.4286
$ws.Range("N116").Value = -7087.3333
$ws.Range("H122").Value = 3272.2856
$ws.Range("I122").Value = 1496
$ws.Range("J122").Value = 3982.8
$ws.Range("K122").Value = 4488
$ws.Range("L122").Value = 11948.4
$ws.Range("M122").Value = -2038
$ws.Range("N122").Value = -16848.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1413.8
$ws.Range("I3").Value = 948.5714
$ws.Range("J3").Value = 2499.3333
$ws.Range("K3").Value = 948.5714
$ws.Range("L3").Value = 2499.3333
$ws.Range("M3").Value = -834.5714
$ws.Range("N3").Value = -2727.3333
$ws.Range("H99").Value = 1450
$ws.Range("I99").Value = 900
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 900
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = 598
$ws.Range("N99").Value = -4996
$ws.Range("H134").Value = 22483.55
$ws.Range("I134").Value = 36408.414
$ws.Range("J134").Value = 2292.5
$ws.Range("K134").Value = 109225.242
$ws.Range("L134").Value = 6877.5
$ws.Range("M134").Value = -106690.242
$ws.Range("N134").Value = -11947.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9092678
$ws.Range("I31").Value = 1828.2778
$ws.Range("J31").Value = 50001500
$ws.Range("K31").Value = 1828.2778
$ws.Range("L31").Value = 50001500
$ws.Range("M31").Value = -1533.2778
$ws.Range("N31").Value = -50002090
$ws.Range("H34").Value = 9092678
$ws.Range("I34").Value = 1828.2778
$ws.Range("J34").Value = 50001500
$ws.Range("K34").Value = 1828.2778
$ws.Range("L34").Value = 50001500
$ws.Range("M34").Value = -1626.2778
$ws.Range("N34").Value = -50001904
$ws.Range("H58").Value = 1090.5927
$ws.Range("I58").Value = 1134
$ws.Range("J58").Value = 987.5
$ws.Range("K58").Value = 1134
$ws.Range("L58").Value = 987.5
$ws.Range("M58").Value = -931
$ws.Range("N58").Value = -1393.5
$ws.Range("H136").Value = 1090.5927
$ws.Range("I136").Value = 1134
$ws.Range("J136").Value = 987.5
$ws.Range("K136").Value = 3402
$ws.Range("L136").Value = 2962.5
$ws.Range("M136").Value = -852
$ws.Range("N136").Value = -8062.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 9043.6
$ws.Range("I33").Value = 10040.111
$ws.Range("J33").Value = 75
$ws.Range("K33").Value = 60240.666
$ws.Range("L33").Value = 450
$ws.Range("M33").Value = -59957.666
$ws.Range("N33").Value = -1016
$ws.Range("H131").Value = 5882.7617
$ws.Range("I131").Value = 14107.5
$ws.Range("J131").Value = 821.38464
$ws.Range("K131").Value = 42322.5
$ws.Range("L131").Value = 2464.15392
$ws.Range("M131").Value = -37282.5
$ws.Range("N131").Value = -12544.15392
$ws.Range("H137").Value = 75784240
$ws.Range("I137").Value = 111114030
$ws.Range("J137").Value = 62535570
$ws.Range("K137").Value = 333342090
$ws.Range("L137").Value = 187606710
$ws.Range("M137").Value = -333336990
$ws.Range("N137").Value = -187616910

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4627.7856
$ws.Range("I80").Value = 2598.75
$ws.Range("J80").Value = 5439.4
$ws.Range("K80").Value = 2598.75
$ws.Range("L80").Value = 5439.4
$ws.Range("M80").Value = -1600.75
$ws.Range("N80").Value = -7435.4
$ws.Range("H83").Value = 4627.7856
$ws.Range("I83").Value = 2598.75
$ws.Range("J83").Value = 5439.4
$ws.Range("K83").Value = 12993.75
$ws.Range("L83").Value = 27197
$ws.Range("M83").Value = -8001.75
$ws.Range("N83").Value = -37181
$ws.Range("H102").Value = 1404.8334
$ws.Range("I102").Value = 1305.8
$ws.Range("J102").Value = 1900
$ws.Range("K102").Value = 1305.8
$ws.Range("L102").Value = 1900
$ws.Range("M102").Value = 316.2
$ws.Range("N102").Value = -5144
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 236.83333
$ws.Range("I55").Value = 267
$ws.Range("J55").Value = 158.4
$ws.Range("K55").Value = 267
$ws.Range("L55").Value = 158.4
$ws.Range("M55").Value = -94
$ws.Range("N55").Value = -504.4
$ws.Range("H132").Value = 6333.275
$ws.Range("I132").Value = 7065.0938
$ws.Range("J132").Value = 3406
$ws.Range("K132").Value = 21195.2814
$ws.Range("L132").Value = 10218
$ws.Range("M132").Value = -18665.2814
$ws.Range("N132").Value = -15278
$ws.Range("H138").Value = 20390
$ws.Range("I138").Value = 20390
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 20390
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -15250

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 330.625
$ws.Range("I107").Value = 315.83334
$ws.Range("J107").Value = 375
$ws.Range("K107").Value = 947.5000200000001
$ws.Range("L107").Value = 1125
$ws.Range("M107").Value = 972.4999799999999
$ws.Range("N107").Value = -4965
$ws.Range("H113").Value = 632
$ws.Range("I113").Value = 720
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 2160
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 10
$ws.Range("N113").Value = -5840
$ws.Range("H122").Value = 1691.3846
$ws.Range("I122").Value = 1936.625
$ws.Range("J122").Value = 1299
$ws.Range("K122").Value = 5809.875
$ws.Range("L122").Value = 3897
$ws.Range("M122").Value = -3359.875
$ws.Range("N122").Value = -8797
$ws.Range("H132").Value = 1121.1136
$ws.Range("I132").Value = 966.9459000000001
$ws.Range("J132").Value = 1936
$ws.Range("K132").Value = 2900.8377
$ws.Range("L132").Value = 5808
$ws.Range("M132").Value = -370.8377
$ws.Range("N132").Value = -10868
